$d = $word.ActiveDocument

function Find-ParaIndexByText($doc, $needle) {
    $n = $doc.Paragraphs.Count
    $found = -1
    for ($i = 1; $i -le $n; $i++) {
        $t = $doc.Paragraphs($i).Range.Text.TrimEnd([char]13)
        if ($t -eq $needle) { $found = $i }
    }
    return $found
}

# ---------------------------------------------------------------------
# Edit 1: After the "General-" paragraph, insert a new paragraph:
#   "When issues " + "are found check for another account before
#   sending to Gary" (two runs), followed by a blank paragraph, before
#   the existing "When sending Multiple chats..." paragraph.
# ---------------------------------------------------------------------
$generalIdx = Find-ParaIndexByText $d "General-"
$generalPara = $d.Paragraphs($generalIdx)
$generalEnd = $generalPara.Range.End
$insPos1 = $generalEnd - 1
$insPoint1 = $d.Range($insPos1, $insPos1)

$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">When issues </w:t></w:r><w:r><w:t>are found check for another account before sending to Gary</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint1.InsertXML($xml1)

# ---------------------------------------------------------------------
# Edit 2: After the "1)Check if multiple accounts are tested" paragraph,
# insert a new paragraph "2)Attach Video to Bug/PBI" (two runs: "2)A"
# and "ttach Video to Bug/PBI"), before the existing blank paragraph.
# ---------------------------------------------------------------------
$checkIdx = Find-ParaIndexByText $d "1)Check if multiple accounts are tested"
$checkPara = $d.Paragraphs($checkIdx)
$checkEnd = $checkPara.Range.End
$insPos2 = $checkEnd - 1
$insPoint2 = $d.Range($insPos2, $insPos2)

$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>2)A</w:t></w:r><w:r><w:t>ttach Video to Bug/PBI</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint2.InsertXML($xml2)

# ---------------------------------------------------------------------
# Edit 3: Remove the "After testing is Done attach Video to Bug/PBI"
# paragraph along with the blank paragraph immediately preceding it
# (the one right after "Before clicking End and Release Device...").
# ---------------------------------------------------------------------
$afterIdx = Find-ParaIndexByText $d "After testing is Done attach Video to Bug/PBI"
$afterPara = $d.Paragraphs($afterIdx)
$afterPara.Range.Delete()

$blankIdx = $afterIdx - 1
$blankPara = $d.Paragraphs($blankIdx)
$blankPara.Range.Delete()

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
